# Apply the "added motif analysis diffbind" change:
# Adds two new data blocks (Homer all peaks / Homer consistent) to the
# "differential peaks" sheet, mirroring the existing "July 2021" /
# "Diffbind Dec 2021" blocks (rows 2:3 and 6:7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("differential peaks")

# --- Block 1: "Homer all peaks" (rows 10:11) ---
$ws.Range("A10").Value = "Less accessible in KO"
$ws.Range("B10").Value = 1658
$ws.Range("C10").Value = "Homer all peaks"

$ws.Range("A11").Value = "More accessible in KO"
$ws.Range("B11").Value = 1879

# --- Block 2: "Homer consistent" (rows 14:15) ---
$ws.Range("A14").Value = "Less accessible in KO"
$ws.Range("B14").Value = 36
$ws.Range("C14").Value = "Homer consistent"

$ws.Range("A15").Value = "More accessible in KO"
$ws.Range("B15").Value = 209

# Copy the formatting of the existing block (rows 6:7), cell by cell, onto
# the two new blocks so the new cells pick up the same styles (borders /
# alignment) already used on the sheet, without creating redundant style
# entries in the workbook's style table.
$ws.Range("A6").Copy(); $ws.Range("A10").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B6").Copy(); $ws.Range("B10").PasteSpecial(-4122)
$ws.Range("C6").Copy(); $ws.Range("C10").PasteSpecial(-4122)
$ws.Range("A7").Copy(); $ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B7").Copy(); $ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C7").Copy(); $ws.Range("C11").PasteSpecial(-4122)

$ws.Range("A6").Copy(); $ws.Range("A14").PasteSpecial(-4122)
$ws.Range("B6").Copy(); $ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C6").Copy(); $ws.Range("C14").PasteSpecial(-4122)
$ws.Range("A7").Copy(); $ws.Range("A15").PasteSpecial(-4122)
$ws.Range("B7").Copy(); $ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C7").Copy(); $ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Merge the label cells in column C, same as the existing blocks.
$ws.Range("C10:C11").Merge()
$ws.Range("C14:C15").Merge()

# Match the final selection recorded in the saved workbook.
$ws.Range("B11").Select()
